# Adds the "(00:28:36)" timestamp, the "Бин JdbcTemplate" section and its
# two highlighted code blocks (the JdbcTemplate bean wiring snippet and the
# spring-data-jdbc doc comment) at the very end of the document, matching the
# "timing 00:41:00 main editing" commit.

$d = $word.ActiveDocument

# Park a fresh, empty paragraph at the end of the story first. That keeps the
# pre-existing final (empty) paragraph completely untouched, in place, ahead
# of the content we are about to insert -- InsertXML always splits whatever
# paragraph the (collapsed) insertion range currently sits inside, and the
# trailing half of that split inherits the old paragraph mark's formatting,
# so we want that split to happen in a brand new paragraph, not in the
# document's existing last paragraph.
$r = $d.Content
$r.Collapse(0)
$r.InsertParagraphAfter() | Out-Null

# Re-collapse to the (new) end of the document and inject the rich OOXML
# fragment describing the new paragraphs, preserving their exact formatting
# (bold/underline heading, "HTML Preformatted" shaded code paragraphs with
# per-token colors, proofErr spell-check markers, manual line breaks, etc.).
$r2 = $d.Content
$r2.Collapse(0)

$xml = @'
<?xml version="1.0" encoding="UTF-8" standalone="yes"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
<w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>(00:28:36)</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/></w:rPr><w:t xml:space="preserve">Бин </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:b/><w:bCs/><w:u w:val="single"/><w:lang w:val="en-US"/></w:rPr><w:t>JdbcTemplate</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:pStyle w:val="HTML"/><w:shd w:val="clear" w:color="auto" w:fill="1E1F22"/><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr></w:pPr><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>ConfigurableApplicationContext</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> context = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>SpringApplication.</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>run</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>SpringDataDemoApplication.</w:t></w:r><w:r><w:rPr><w:color w:val="CF8E6D"/><w:lang w:val="en-US"/></w:rPr><w:t>class</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>args</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>);</w:t></w:r><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>JdbcTemplate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>jdbcTemplate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>context.getBean</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>JdbcTemplate.</w:t></w:r><w:r><w:rPr><w:color w:val="CF8E6D"/><w:lang w:val="en-US"/></w:rPr><w:t>class</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>);</w:t></w:r><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:color w:val="7A7E85"/><w:lang w:val="en-US"/></w:rPr><w:t>//</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="7A7E85"/><w:lang w:val="en-US"/></w:rPr><w:t>JdbcTemplate</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="7A7E85"/><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:color w:val="7A7E85"/><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>jdbcTemplate.execute</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:color w:val="6AAB73"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">"create table users(id </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="6AAB73"/><w:lang w:val="en-US"/></w:rPr><w:t>bigint</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="6AAB73"/><w:lang w:val="en-US"/></w:rPr><w:t>, name varchar(512))"</w:t></w:r><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>);</w:t></w:r><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>jdbcTemplate.execute</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:color w:val="6AAB73"/><w:lang w:val="en-US"/></w:rPr><w:t>"insert into users(id, name) values(1, 'Igor')"</w:t></w:r><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>);</w:t></w:r><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:br/><w:t xml:space="preserve">List&lt;User&gt; users = </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>jdbcTemplate.query</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:color w:val="6AAB73"/><w:lang w:val="en-US"/></w:rPr><w:t>"select id, name from users"</w:t></w:r><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:color w:val="CF8E6D"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">new </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>RowMapper</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>&lt;User&gt;() {</w:t></w:r><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:br/><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr><w:color w:val="B3AE60"/><w:lang w:val="en-US"/></w:rPr><w:t>@Override</w:t></w:r><w:r><w:rPr><w:color w:val="B3AE60"/><w:lang w:val="en-US"/></w:rPr><w:br/><w:t xml:space="preserve">    </w:t></w:r><w:r><w:rPr><w:color w:val="CF8E6D"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">public </w:t></w:r><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">User </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="56A8F5"/><w:lang w:val="en-US"/></w:rPr><w:t>mapRow</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>ResultSet</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>rs</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">, </w:t></w:r><w:r><w:rPr><w:color w:val="CF8E6D"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">int </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>rowNum</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">) </w:t></w:r><w:r><w:rPr><w:color w:val="CF8E6D"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">throws </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>SQLException</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve"> {</w:t></w:r><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:br/><w:t xml:space="preserve">       </w:t></w:r><w:r><w:rPr><w:color w:val="CF8E6D"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">return new </w:t></w:r><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>User(</w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>rs.getLong</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:color w:val="6AAB73"/><w:lang w:val="en-US"/></w:rPr><w:t>"id"</w:t></w:r><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t xml:space="preserve">), </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>rs.getString</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>(</w:t></w:r><w:r><w:rPr><w:color w:val="6AAB73"/><w:lang w:val="en-US"/></w:rPr><w:t>"name"</w:t></w:r><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:t>));</w:t></w:r><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:br/><w:t xml:space="preserve">    }</w:t></w:r><w:r><w:rPr><w:color w:val="BCBEC4"/><w:lang w:val="en-US"/></w:rPr><w:br/><w:t>});</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="HTML"/><w:shd w:val="clear" w:color="auto" w:fill="1E1F22"/><w:rPr><w:color w:val="BCBEC4"/></w:rPr></w:pPr><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="5F826B"/></w:rPr><w:t>/**</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="5F826B"/></w:rPr><w:br/><w:t xml:space="preserve"> * </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="5F826B"/></w:rPr><w:t>spring-data-jdbc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="5F826B"/></w:rPr><w:t xml:space="preserve"> - набор готовых инструментов для взаимодействия с базой данных.</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="5F826B"/></w:rPr><w:br/></w:r><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="5F826B"/></w:rPr><w:lastRenderedPageBreak/><w:t xml:space="preserve"> </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="5F826B"/></w:rPr><w:t>* По сути</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="5F826B"/></w:rPr><w:t xml:space="preserve">, оборачивает стандартный JDBC и предоставляет удобные интерфейсы для настройки и </w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="5F826B"/></w:rPr><w:br/><w:t xml:space="preserve"> * взаимодействия с БД.</w:t></w:r><w:r><w:rPr><w:i/><w:iCs/><w:color w:val="5F826B"/></w:rPr><w:br/><w:t xml:space="preserve"> */</w:t></w:r></w:p><w:p/><w:p/>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$r2.InsertXML($xml) | Out-Null

# InsertXML splits the paragraph the range was collapsed inside (the blank
# paragraph opened above) and carries that paragraph's original mark/rPr to
# the very end of the document as a surplus empty paragraph after the
# fragment is spliced in. Delete it so the document ends with exactly the two
# bare empty paragraphs that are part of the inserted fragment itself.
$d.Paragraphs.Last.Range.Delete() | Out-Null

Write-Output ("Paragraphs: " + $d.Paragraphs.Count)
